$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 updates
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9

# Row 14 updates
$ws.Range("H14").Value = 3.25
$ws.Range("I14").Value = 2.52
$ws.Range("K14").Value = 2.12
$ws.Range("L14").Value = 3.1
$ws.Range("N14").Value = 7
$ws.Range("O14").Value = 1.34
$ws.Range("P14").Value = 3
$ws.Range("Q14").Value = 2.02
$ws.Range("U14").Value = 1.8
$ws.Range("V14").Value = 1.91
$ws.Range("AC14").Value = 7
$ws.Range("AD14").Value = 6.4
$ws.Range("AE14").Value = 14.5
$ws.Range("AH14").Value = 8
$ws.Range("AK14").Value = 27
$ws.Range("AO14").Value = 13.5
$ws.Range("AU14").Value = 7.2
$ws.Range("AW14").Value = 4.45
$ws.Range("AY14").Value = 22
$ws.Range("BB14").Value = 300
